$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 94
$ws.Range("H94").Value = 2900.7144
$ws.Range("I94").Value = 2676.25
$ws.Range("J94").Value = 3200
$ws.Range("K94").Value = 2676.25
$ws.Range("L94").Value = 3200
$ws.Range("M94").Value = -2225.25
$ws.Range("N94").Value = -4102
# Row 113
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3502
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3502
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -10010
# Row 127
$ws.Range("H127").Value = 1588.2
$ws.Range("I127").Value = 924.6667
$ws.Range("J127").Value = 2583.5
$ws.Range("K127").Value = 2774.0001
$ws.Range("L127").Value = 7750.5
$ws.Range("M127").Value = 2185.9999
$ws.Range("N127").Value = -17670.5
# Row 137
$ws.Range("H137").Value = 1248.5217
$ws.Range("I137").Value = 1116.2051
$ws.Range("K137").Value = 3348.615299999999
$ws.Range("M137").Value = -798.6152999999995

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 43
$ws.Range("H43").Value = 7897
$ws.Range("J43").Value = 7897
$ws.Range("L43").Value = 7897
$ws.Range("N43").Value = -8523
# Row 102
$ws.Range("H102").Value = 1683.3334
$ws.Range("I102").Value = 1683.3334
$ws.Range("K102").Value = 1683.3334
$ws.Range("M102").Value = -61.33339999999998
# Row 104
$ws.Range("H104").Value = 42390
$ws.Range("J104").Value = 42390
$ws.Range("L104").Value = 42390
$ws.Range("N104").Value = -49378
# Row 122
$ws.Range("H122").Value = 1727.15
$ws.Range("I122").Value = 1196
$ws.Range("J122").Value = 2713.5715
$ws.Range("K122").Value = 3588
$ws.Range("L122").Value = 8140.7145
$ws.Range("M122").Value = -1138
$ws.Range("N122").Value = -13040.7145

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1105.0555
$ws.Range("I94").Value = 846.4
$ws.Range("J94").Value = 2398.3333
$ws.Range("K94").Value = 846.4
$ws.Range("L94").Value = 2398.3333
$ws.Range("M94").Value = -395.4
$ws.Range("N94").Value = -3300.3333
# Row 99
$ws.Range("H99").Value = 811.5
$ws.Range("I99").Value = 573.8
$ws.Range("K99").Value = 573.8
$ws.Range("M99").Value = 924.2
# Row 103
$ws.Range("H103").Value = 32500
$ws.Range("J103").Value = 32500
$ws.Range("L103").Value = 32500
$ws.Range("N103").Value = -34844
# Row 105
$ws.Range("H105").Value = 2377.5
$ws.Range("I105").Value = 2336.6667
$ws.Range("K105").Value = 2336.6667
$ws.Range("M105").Value = -589.6667000000002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1470.1538
$ws.Range("I58").Value = 1468
$ws.Range("J58").Value = 1475
$ws.Range("K58").Value = 1468
$ws.Range("L58").Value = 1475
$ws.Range("M58").Value = -1265
$ws.Range("N58").Value = -1881
# Row 94
$ws.Range("H94").Value = 55556840
$ws.Range("I94").Value = 142857620
$ws.Range("J94").Value = 1801.8182
$ws.Range("K94").Value = 142857620
$ws.Range("L94").Value = 1801.8182
$ws.Range("M94").Value = -142857169
$ws.Range("N94").Value = -2703.8182
# Row 99
$ws.Range("H99").Value = 2538.4614
$ws.Range("I99").Value = 1900
$ws.Range("J99").Value = 2937.5
$ws.Range("K99").Value = 1900
$ws.Range("L99").Value = 2937.5
$ws.Range("M99").Value = -402
$ws.Range("N99").Value = -5933.5
# Row 126
$ws.Range("H126").Value = 2538.4614
$ws.Range("I126").Value = 1900
$ws.Range("J126").Value = 2937.5
$ws.Range("K126").Value = 5700
$ws.Range("L126").Value = 8812.5
$ws.Range("M126").Value = -3230
$ws.Range("N126").Value = -13752.5
# Row 132
$ws.Range("H132").Value = 3777.5557
$ws.Range("I132").Value = 3076.7693
$ws.Range("J132").Value = 5599.6
$ws.Range("K132").Value = 9230.3079
$ws.Range("L132").Value = 16798.8
$ws.Range("M132").Value = -6700.3079
$ws.Range("N132").Value = -21858.8
# Row 134
$ws.Range("H134").Value = 1016.8461
$ws.Range("I134").Value = 997.52
$ws.Range("K134").Value = 2992.56
$ws.Range("M134").Value = -457.5599999999999
# Row 135
$ws.Range("H135").Value = 35000
$ws.Range("J135").Value = 35000
$ws.Range("L135").Value = 35000
$ws.Range("N135").Value = -45140
# Row 136
$ws.Range("H136").Value = 1470.1538
$ws.Range("I136").Value = 1468
$ws.Range("J136").Value = 1475
$ws.Range("K136").Value = 4404
$ws.Range("L136").Value = 4425
$ws.Range("M136").Value = -1854
$ws.Range("N136").Value = -9525

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 98
$ws.Range("H98").Value = 32500.5
$ws.Range("J98").Value = 32500.5
$ws.Range("L98").Value = 32500.5
$ws.Range("N98").Value = -38490.5
# Row 102
$ws.Range("H102").Value = 1750
$ws.Range("I102").Value = 1500
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1500
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 122
$ws.Range("N102").Value = -5244
# Row 126
$ws.Range("H126").Value = 6668389
$ws.Range("J126").Value = 8334986
$ws.Range("L126").Value = 25004958
$ws.Range("N126").Value = -25009898

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2580.6667
$ws.Range("I40").Value = 2580.6667
$ws.Range("K40").Value = 2580.6667
$ws.Range("M40").Value = -2444.6667
# Row 93
$ws.Range("H93").Value = 1423537.6
$ws.Range("I93").Value = 2080008.9
$ws.Range("J93").Value = 1183
$ws.Range("K93").Value = 2080008.9
$ws.Range("L93").Value = 1183
$ws.Range("M93").Value = -2078760.9
$ws.Range("N93").Value = -3679
# Row 122
$ws.Range("H122").Value = 3225.2
$ws.Range("I122").Value = 3307.0908
$ws.Range("K122").Value = 9921.2724
$ws.Range("M122").Value = -7471.2724

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 76
$ws.Range("H76").Value = 19862.166
$ws.Range("J76").Value = 20724.334
$ws.Range("L76").Value = 20724.334
$ws.Range("N76").Value = -21354.334
# Row 79
$ws.Range("H79").Value = 19862.166
$ws.Range("J79").Value = 20724.334
$ws.Range("L79").Value = 20724.334
$ws.Range("N79").Value = -22908.334
# Row 122
$ws.Range("H122").Value = 2429.2856
$ws.Range("I122").Value = 2700
$ws.Range("J122").Value = 1752.5
$ws.Range("K122").Value = 8100
$ws.Range("L122").Value = 5257.5
$ws.Range("M122").Value = -5650
$ws.Range("N122").Value = -10157.5
# Row 126
$ws.Range("H126").Value = 100004
$ws.Range("I126").Value = 100004
$ws.Range("K126").Value = 300012
$ws.Range("M126").Value = -297542
